$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 29493.645
$ws.Range("J17").Value = 30436.768
$ws.Range("L17").Value = 91310.304
$ws.Range("N17").Value = -91646.304

# Row 43
$ws.Range("H43").Value = 6144.8125
$ws.Range("I43").Value = 5939.8335
$ws.Range("J43").Value = 6267.8
$ws.Range("K43").Value = 5939.8335
$ws.Range("L43").Value = 6267.8
$ws.Range("M43").Value = -5870.8335
$ws.Range("N43").Value = -6405.8

# Row 70
$ws.Range("H70").Value = 1087.6666
$ws.Range("I70").Value = 896.3333
$ws.Range("J70").Value = 1183.3334
$ws.Range("K70").Value = 2688.9999
$ws.Range("L70").Value = 3550.0002
$ws.Range("M70").Value = -2418.9999
$ws.Range("N70").Value = -4090.0002

# Row 73
$ws.Range("H73").Value = 1087.6666
$ws.Range("I73").Value = 896.3333
$ws.Range("J73").Value = 1183.3334
$ws.Range("K73").Value = 2688.9999
$ws.Range("L73").Value = 3550.0002
$ws.Range("M73").Value = -1752.9999
$ws.Range("N73").Value = -5422.0002

# Row 80
$ws.Range("H80").Value = 603.1
$ws.Range("I80").Value = 446.66666
$ws.Range("J80").Value = 670.1429000000001
$ws.Range("K80").Value = 1339.99998
$ws.Range("L80").Value = 2010.4287
$ws.Range("M80").Value = -341.9999800000001
$ws.Range("N80").Value = -4006.4287

# Row 83
$ws.Range("H83").Value = 603.1
$ws.Range("I83").Value = 446.66666
$ws.Range("J83").Value = 670.1429000000001
$ws.Range("K83").Value = 4019.99994
$ws.Range("L83").Value = 6031.2861
$ws.Range("M83").Value = 972.0000600000003
$ws.Range("N83").Value = -16015.2861

# Row 92
$ws.Range("H92").Value = 1304
$ws.Range("I92").Value = 1304
$ws.Range("K92").Value = 1304
$ws.Range("M92").Value = -56

# Row 116
$ws.Range("H116").Value = 5105.5557
$ws.Range("I116").Value = 5316.6665
$ws.Range("K116").Value = 5316.6665
$ws.Range("M116").Value = -1874.6665

# Row 137
$ws.Range("H137").Value = 3275.3958
$ws.Range("I137").Value = 1283.5555
$ws.Range("K137").Value = 3850.6665
$ws.Range("M137").Value = -1300.6665


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 1224.7368
$ws.Range("I61").Value = 1027.6471
$ws.Range("J61").Value = 2900
$ws.Range("K61").Value = 1027.6471
$ws.Range("L61").Value = 2900
$ws.Range("M61").Value = -815.6470999999999
$ws.Range("N61").Value = -3324

# Row 74
$ws.Range("H74").Value = 2160.561
$ws.Range("I74").Value = 998.3077
$ws.Range("K74").Value = 998.3077
$ws.Range("M74").Value = -124.3077

# Row 77
$ws.Range("H77").Value = 2160.561
$ws.Range("I77").Value = 998.3077
$ws.Range("K77").Value = 4991.5385
$ws.Range("M77").Value = -623.5384999999997

# Row 122
$ws.Range("H122").Value = 3902.4443
$ws.Range("I122").Value = 2437
$ws.Range("K122").Value = 7311
$ws.Range("M122").Value = -4861

# Row 132
$ws.Range("H132").Value = 15939.286
$ws.Range("I132").Value = 17937.5
$ws.Range("K132").Value = 53812.5
$ws.Range("M132").Value = -51282.5

# Row 136
$ws.Range("H136").Value = 1224.7368
$ws.Range("I136").Value = 1027.6471
$ws.Range("J136").Value = 2900
$ws.Range("K136").Value = 3082.9413
$ws.Range("L136").Value = 8700
$ws.Range("M136").Value = -532.9412999999995
$ws.Range("N136").Value = -13800


# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 75
$ws.Range("H75").Value = 15166.667
$ws.Range("I75").Value = 15166.667
$ws.Range("K75").Value = 15166.667
$ws.Range("M75").Value = -14230.667

# Row 76
$ws.Range("H76").Value = 1699.5
$ws.Range("J76").Value = 1699.5
$ws.Range("L76").Value = 1699.5
$ws.Range("N76").Value = -2329.5

# Row 78
$ws.Range("H78").Value = 15166.667
$ws.Range("I78").Value = 15166.667
$ws.Range("K78").Value = 45500.001
$ws.Range("M78").Value = -40820.001

# Row 79
$ws.Range("H79").Value = 1699.5
$ws.Range("J79").Value = 1699.5
$ws.Range("L79").Value = 1699.5
$ws.Range("N79").Value = -3883.5


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2427.3914
$ws.Range("I31").Value = 1844.091
$ws.Range("K31").Value = 1844.091
$ws.Range("M31").Value = -1549.091

# Row 34
$ws.Range("H34").Value = 2427.3914
$ws.Range("I34").Value = 1844.091
$ws.Range("K34").Value = 1844.091
$ws.Range("M34").Value = -1642.091

# Row 132
$ws.Range("H132").Value = 5057.375
$ws.Range("I132").Value = 5824.75
$ws.Range("J132").Value = 2755.25
$ws.Range("K132").Value = 17474.25
$ws.Range("L132").Value = 8265.75
$ws.Range("M132").Value = -14944.25
$ws.Range("N132").Value = -13325.75


# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 137
$ws.Range("H137").Value = 4549165
$ws.Range("I137").Value = 10001516
$ws.Range("J137").Value = 5539.0835
$ws.Range("K137").Value = 30004548
$ws.Range("L137").Value = 16617.2505
$ws.Range("M137").Value = -29999448
$ws.Range("N137").Value = -26817.2505


# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 4651
$ws.Range("I132").Value = 4597.6665
$ws.Range("K132").Value = 13792.9995
$ws.Range("M132").Value = -11262.9995


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 233.66667
$ws.Range("I22").Value = 225.5
$ws.Range("J22").Value = 250
$ws.Range("K22").Value = 225.5
$ws.Range("L22").Value = 250
$ws.Range("M22").Value = 69.5
$ws.Range("N22").Value = -840

# Row 27
$ws.Range("H27").Value = 233.66667
$ws.Range("I27").Value = 225.5
$ws.Range("J27").Value = 250
$ws.Range("K27").Value = 225.5
$ws.Range("L27").Value = 250
$ws.Range("M27").Value = -118.5
$ws.Range("N27").Value = -464

# Row 93
$ws.Range("H93").Value = 35814.5
$ws.Range("I93").Value = 2602.4285
$ws.Range("J93").Value = 113309.336
$ws.Range("K93").Value = 2602.4285
$ws.Range("L93").Value = 113309.336
$ws.Range("M93").Value = -1354.4285
$ws.Range("N93").Value = -115805.336

# Row 132
$ws.Range("H132").Value = 2809.1904
$ws.Range("I132").Value = 2549.9
$ws.Range("K132").Value = 7649.700000000001
$ws.Range("M132").Value = -5119.700000000001

# Row 136
$ws.Range("H136").Value = 3669.818
$ws.Range("I136").Value = 3086.2
$ws.Range("K136").Value = 9258.599999999999
$ws.Range("M136").Value = -6708.599999999999


# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 3179.4546
$ws.Range("I81").Value = 3714.647
$ws.Range("J81").Value = 1359.8
$ws.Range("K81").Value = 7429.294
$ws.Range("L81").Value = 2719.6
$ws.Range("M81").Value = -6368.294
$ws.Range("N81").Value = -4841.6

# Row 84
$ws.Range("H84").Value = 3179.4546
$ws.Range("I84").Value = 3714.647
$ws.Range("J84").Value = 1359.8
$ws.Range("K84").Value = 37146.47
$ws.Range("L84").Value = 13598
$ws.Range("M84").Value = -31842.47
$ws.Range("N84").Value = -24206

# Row 132
$ws.Range("H132").Value = 17599.209
$ws.Range("I132").Value = 16814.64
$ws.Range("J132").Value = 25248.75
$ws.Range("K132").Value = 50443.92
$ws.Range("L132").Value = 75746.25
$ws.Range("M132").Value = -47913.92
$ws.Range("N132").Value = -80806.25

# Row 136
$ws.Range("H136").Value = 2120.7441
$ws.Range("I136").Value = 2183.7297
$ws.Range("K136").Value = 6551.1891
$ws.Range("M136").Value = -4001.1891

# Row 139
$ws.Range("H139").Value = 99995
$ws.Range("I139").Value = 99995
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 99995
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -94855
$ws.Range("N139").ClearContents()

